$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Update Shadow Rate with Latest Data" ---
# Re-run of the shadow-rate estimation against the latest vintage of the
# underlying fed funds rate series. The refresh nudges the cached
# "fedfundsrate" input for 1984-Q1 by a hair (B2), rewrites the modeled
# "fedfundsrate_shadow" series (column C) for every historical quarter, and
# appends five newly observed quarters (2020-Q1 .. 2021-Q1) during which the
# funds rate sits at the zero lower bound while the shadow rate goes negative.

# Tiny data revision to the 1984-Q1 fed funds rate input
$ws.Range("B2").Value = 9.6866666666665768

# Recomputed fedfundsrate_shadow (column C) for existing rows 2-128
$shadowRateRows = @(
    2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 21, 22, 23, 24, 25, 26, 30, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 110, 111, 112, 113, 114, 115, 116, 117, 118, 119, 120, 121, 122, 123, 124, 125, 126, 127, 128
)
$shadowRateValues = @(
    9.6866666666665768, 10.556666666666569, 11.389999999999878, 9.2666666666666231, 8.4766666666667092, 7.9233333333333045, 7.9000000000000181, 8.1033333333333513, 7.8266666666665596, 6.9199999999998596, 6.2066666666666048, 6.266666666666687, 6.2200000000000033, 6.6499999999999782, 6.843333333333379, 6.9166666666666377, 6.6633333333333544, 7.1566666666666556, 8.4699999999999331, 9.4433333333332925, 9.7266666666665724, 9.0833333333331989, 8.6133333333332516, 8.2499999999999574, 6.426666666666625, 1.7400000000000304, 1.4433333333333298, 1.2499999999999956, 1.246666666666707, 1.0166666666666879, 0.99638312099048498, 1.0029047998475038, 1.0095083296175922, 1.4328298615421842, 1.949519841192604, 2.4695666422738238, 2.9429614988056185, 3.4596973846305001, 3.9797693220471109, 4.4565068918595463, 4.9065739074908077, 5.2466350660383743, 5.2466889058359056, 5.2567344970156249, 5.2501048540359641, 5.073466484556155, 4.4968192910172844, 3.1768298489961877, 2.0868331992731415, 1.9401643070454666, 0.50777465922533516, 1.6488952284639868, 0.62422887211277533, -0.11067726072560014, -0.16189799346500999, -0.12751431539200775, -1.1848463498480011, -1.0859586593229498, -1.5402196670599944, -1.3623043248304323, -1.0891450601857144, -1.8674739911055394, -1.2162518929261901, -2.009660476457209, -1.6607243994285259, -1.2895588468961661, -2.4213785264856846, -1.1357874445460658, -0.63276879257986307, -0.45619773641784045, -0.54336708914966003, -0.74278953819548743, -0.72762278155045212, -0.4624926551333397, -0.21961124805635057, 0.16225716019162562, 0.032963145757958223, 0.05031697791864076
)
for ($i = 0; $i -lt $shadowRateRows.Count; $i++) {
    $ws.Cells.Item($shadowRateRows[$i], 3).Value = $shadowRateValues[$i]
}

# Newly observed quarters appended to the series (2020-Q1 .. 2021-Q1)
$newDates    = @(2020, 2020.25, 2020.5, 2020.75, 2021)
$newFunds    = @(1.2599999999999723, 0, 0, 0, 0)
$newShadow   = @(1.2599999999999723, 7.6631433313967046, -4.7850751537857779, -3.048671875166753, -3.4569946319478295)
$startRow = 146
for ($i = 0; $i -lt $newDates.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = $newFunds[$i]
    $ws.Cells.Item($r, 3).Value = $newShadow[$i]
}
